$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks first (row deletion below does not keep the
# hyperlink collection in sync), then drop the three stale project rows
# (4, 5, 6) so only the header and the two refreshed listings remain.
$ws.Cells.Hyperlinks.Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Refresh row 2 with the newly scraped listing.
$ws.Range("A2").Value = "2025-12-31 06:30:30"
$ws.Range("B2").Value = "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5450864"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5450864") | Out-Null
$ws.Range("G2").Value = 383
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# Refresh row 3 with the newly scraped listing.
$ws.Range("A3").Value = "2025-12-31 06:30:30"
$ws.Range("B3").Value = "実績づくり歓迎。既存のものをベースに、自動化や軽い修正をお願いできる方を探しています"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5463636"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5463636") | Out-Null
$ws.Range("G3").Value = 80
$ws.Range("H3").Value = "◆自動化"

# Column B narrows by one character unit (52 -> 51). ColumnWidth round-trips
# through a pixel conversion, so 50.1667 is the input that lands on an
# XML-serialized width of exactly 51.
$ws.Columns.Item(2).ColumnWidth = 50.1667
